# Updated cryptos list on Wed Jun 14 18:48:54 UTC 2023 with GitHub Actions
# (refreshed Price / Volume(1h) snapshot from coinranking.com, plus a couple
# of rows that swapped ranking order: TRON/WrappedEther and Toncoin/LidoDAOToken)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.926.65'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.736.50'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9987'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.76%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9996'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5032'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2725'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06178'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.89%  '
$ws.Range('B10').Value = 'WrappedEther'
$ws.Range('C10').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.742.19'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.19%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07252'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.07%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.6571'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.20'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.778'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.17'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.06%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9999'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9981'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.921.89'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.86'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000006819'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.96%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.596'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +7.98%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.960.43'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.798'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.57%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.478'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.73%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '133.85'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.43%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.23'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.60%  '
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.423'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.70%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.790'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '105.68'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.994'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08123'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.716'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04729'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.20%  '
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.002'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6124'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.78%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.743'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01604'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.8587'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +16.81%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.955'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9992'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '100.61'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.21%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.3918'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.022'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.56%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1180'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.97%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.328'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.95%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.86'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.97%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05275'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.39%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.74'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.3484'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.35%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.629'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.95%  '
